# Updated symbol list on Wed Dec 21 18:55:33 UTC 2022 with GitHub Actions
# Applies the refreshed price / symbol snapshot to the cryptos sheet.
#
# Numeric-looking values are written with a leading apostrophe so Excel
# keeps them as literal text (matching the workbook's inline-string cells)
# instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.69"
$ws.Range("D3").Value = "'22.42"
$ws.Range("D4").Value = "'5.248"
$ws.Range("D5").Value = "'0.05684"
$ws.Range("D7").Value = "'6.310"
$ws.Range("D8").Value = "'0.8065"
$ws.Range("D9").Value = "'0.8694"
$ws.Range("D10").Value = "'0.1414"
$ws.Range("D11").Value = "'0.07404"
$ws.Range("D12").Value = "'0.03048"
$ws.Range("D13").Value = "'0.03077"
$ws.Range("D14").Value = "'0.09387"
$ws.Range("D15").Value = "'3.868"
$ws.Range("D16").Value = "'0.001575"
$ws.Range("D17").Value = "'0.04776"

# Rows 18-20 were re-ranked: UpBots -> One -> TigerCash -> UpBots shuffle.
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005812"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006402"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "UpBots"
$ws.Range("C20").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D20").Value = "'0.008424"
$ws.Range("E20").Value = "19UpBotsUBXTBestin24h"

$ws.Range("D21").Value = "'0.005034"
$ws.Range("D22").Value = "'0.0009964"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("D24").Value = "'3.692"
$ws.Range("D25").Value = "'2.193"
$ws.Range("D26").Value = "'0.3277"
$ws.Range("D27").Value = "'0.1303"

$ws.Range("D41").Value = "'0.006851"
$ws.Range("D42").Value = "'0.1064"
$ws.Range("D43").Value = "'0.003201"
$ws.Range("D44").Value = "'0.008452"
$ws.Range("D45").Value = "'0.00005593"

$ws.Range("D48").Value = "'0.2023"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$ws.Range("D50").Value = "'0.01010"
